# Updates cryptos list values (Price and Volume(1h) columns) per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.239.08'
$ws.Range("E2").Value = '  -0.93%  '
$ws.Range("D3").Value = '1.661.59'
$ws.Range("E3").Value = '  -0.96%  '
$ws.Range("D4").Value = '''1.004'
$ws.Range("E4").Value = '  +0.32%  '
$ws.Range("D5").Value = '''218.93'
$ws.Range("E5").Value = '  +1.14%  '
$ws.Range("D6").Value = '''0.5224'
$ws.Range("E6").Value = '  -1.76%  '
$ws.Range("E7").Value = '  +0.36%  '
$ws.Range("D8").Value = '''0.2672'
$ws.Range("E8").Value = '  -0.07%  '
$ws.Range("D9").Value = '''0.06340'
$ws.Range("E9").Value = '  -0.82%  '
$ws.Range("E10").Value = '  -2.51%  '
$ws.Range("D11").Value = '''0.07709'
$ws.Range("E11").Value = '  -1.18%  '
$ws.Range("D12").Value = '1.665.65'
$ws.Range("E12").Value = '  -0.82%  '
$ws.Range("D13").Value = '''4.431'
$ws.Range("E13").Value = '  -1.48%  '
$ws.Range("D14").Value = '1.890.41'
$ws.Range("E14").Value = '  -0.85%  '
$ws.Range("D15").Value = '''0.5475'
$ws.Range("E15").Value = '  -1.59%  '
$ws.Range("D16").Value = '0.0₅8216'
$ws.Range("E16").Value = '  -1.62%  '
$ws.Range("D17").Value = '''65.02'
$ws.Range("E17").Value = '  -0.95%  '
$ws.Range("D18").Value = '26.275.85'
$ws.Range("E18").Value = '  -0.92%  '
$ws.Range("E19").Value = '  +0.43%  '
$ws.Range("D20").Value = '''4.656'
$ws.Range("E20").Value = '  -2.18%  '
$ws.Range("D21").Value = '''195.48'
$ws.Range("E21").Value = '  +0.42%  '
$ws.Range("E22").Value = '  -2.14%  '
$ws.Range("D23").Value = '''6.088'
$ws.Range("E23").Value = '  -3.90%  '
$ws.Range("E24").Value = '  +0.56%  '
$ws.Range("D25").Value = '''138.82'
$ws.Range("E25").Value = '  -3.47%  '
$ws.Range("D26").Value = '''0.1239'
$ws.Range("E26").Value = '  -3.18%  '
$ws.Range("D27").Value = '''7.245'
$ws.Range("E27").Value = '  -2.51%  '
$ws.Range("D28").Value = '''16.18'
$ws.Range("E28").Value = '  -0.84%  '
$ws.Range("E29").Value = '  -0.67%  '
$ws.Range("D30").Value = '''0.05949'
$ws.Range("E30").Value = '  -3.21%  '
$ws.Range("D32").Value = '''3.635'
$ws.Range("E32").Value = '  +0.47%  '
$ws.Range("D33").Value = '''3.308'
$ws.Range("E33").Value = '  -4.11%  '
$ws.Range("D34").Value = '''1.632'
$ws.Range("E34").Value = '  -3.45%  '
$ws.Range("D35").Value = '''0.9795'
$ws.Range("E35").Value = '  -2.62%  '
$ws.Range("D36").Value = '''2.420'
$ws.Range("E36").Value = '  -0.30%  '
$ws.Range("E37").Value = '  +0.02%  '
$ws.Range("D38").Value = '''0.5892'
$ws.Range("E38").Value = '  +2.60%  '
$ws.Range("D39").Value = '''0.01595'
$ws.Range("E39").Value = '  -2.57%  '
$ws.Range("D40").Value = '''5.979'
$ws.Range("E40").Value = '  -0.91%  '
$ws.Range("D41").Value = '''0.8596'
$ws.Range("E41").Value = '  -0.10%  '
$ws.Range("D43").Value = '1.028.31'
$ws.Range("E43").Value = '  -4.14%  '
$ws.Range("D44").Value = '''99.75'
$ws.Range("E44").Value = '  -0.27%  '
$ws.Range("E45").Value = '  -1.19%  '
$ws.Range("E46").Value = '  +7.39%  '
$ws.Range("D47").Value = '''57.33'
$ws.Range("E47").Value = '  +0.52%  '
$ws.Range("D48").Value = '''1.006'
$ws.Range("E48").Value = '  +0.40%  '
$ws.Range("D49").Value = '''8.092'
$ws.Range("E49").Value = '  -0.56%  '
$ws.Range("D50").Value = '''0.05186'
$ws.Range("E50").Value = '  -0.38%  '
$ws.Range("D51").Value = '''1.467'
$ws.Range("E51").Value = '  +0.14%  '
